$d = $word.ActiveDocument

# The document currently ends with two empty paragraphs right before the
# section break. We leave the first one alone, turn the last (empty)
# paragraph into the new "calculator app" project line, and add a brand
# new paragraph after it containing a hyperlink to the demo video
# (mirroring the existing project entries above it).

$count = $d.Paragraphs.Count
$textParagraph = $d.Paragraphs.Item($count)
$textParagraph.Range.Text = "calculator app JavaScript"

# Add a brand-new paragraph after it to host the hyperlink + trailing space,
# just like the other "project link" paragraphs in this document.
$textParagraph = $d.Paragraphs.Item($count)
$insertionPoint = $d.Range($textParagraph.Range.End, $textParagraph.Range.End)
$insertionPoint.InsertParagraphAfter()

$linkParagraph = $d.Paragraphs.Item($count + 1)
$url = "https://www.youtube.com/watch?v=MaV9VqFMzB4&ab_channel=ProgHub"
$linkParagraph.Range.Text = $url + " "

$linkParagraph = $d.Paragraphs.Item($count + 1)
$urlRange = $d.Range($linkParagraph.Range.Start, $linkParagraph.Range.Start + $url.Length)
$d.Hyperlinks.Add($urlRange, $url, "", "", $url)
